$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.12221348285675
$ws.Range("B1").Value = 2.415584087371826
$ws.Range("C1").Value = 5.149642467498779
$ws.Range("D1").Value = 2.258928775787354
$ws.Range("E1").Value = 1.265967965126038
